$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$table.Cell(1, 1).Range.Text = "11÷7=1, 4"
$table.Cell(1, 2).Range.Text = "36÷6=6, 0"
$table.Cell(1, 3).Range.Text = "94÷7=13, 3"
$table.Cell(1, 4).Range.Text = "65÷3=21, 2"
$table.Cell(1, 5).Range.Text = "58÷7=8, 2"
$table.Cell(5, 1).Range.Text = "50÷4=12, 2"
$table.Cell(5, 2).Range.Text = "87÷2=43, 1"
$table.Cell(5, 3).Range.Text = "92÷2=46, 0"
$table.Cell(5, 4).Range.Text = "27÷6=4, 3"
$table.Cell(5, 5).Range.Text = "44÷2=22, 0"
$table.Cell(9, 1).Range.Text = "29÷5=5, 4"
$table.Cell(9, 2).Range.Text = "86÷3=28, 2"
$table.Cell(9, 3).Range.Text = "33÷3=11, 0"
$table.Cell(9, 4).Range.Text = "76÷3=25, 1"
$table.Cell(9, 5).Range.Text = "43÷5=8, 3"
$table.Cell(13, 1).Range.Text = "96÷2=48, 0"
$table.Cell(13, 2).Range.Text = "45÷5=9, 0"
$table.Cell(13, 3).Range.Text = "46÷4=11, 2"
$table.Cell(13, 4).Range.Text = "40÷9=4, 4"
$table.Cell(13, 5).Range.Text = "89÷8=11, 1"
$table.Cell(17, 1).Range.Text = "86÷2=43, 0"
$table.Cell(17, 2).Range.Text = "31÷6=5, 1"
$table.Cell(17, 3).Range.Text = "89÷3=29, 2"
$table.Cell(17, 4).Range.Text = "49÷8=6, 1"
$table.Cell(17, 5).Range.Text = "17÷2=8, 1"
